$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Benchmark summary values (top of the results table)
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "2003"
$t.Cell(5,1).Range.Text  = "0.00002"
$t.Cell(6,1).Range.Text  = "0.58457"
$t.Cell(7,1).Range.Text  = "0.07722"
$t.Cell(8,1).Range.Text  = "0.02883"
$t.Cell(9,1).Range.Text  = "0.34051"
$t.Cell(10,1).Range.Text = "0.34051"
$t.Cell(11,1).Range.Text = "0.58457"
$t.Cell(12,1).Range.Text = "1.29695"

# Final three rows: collapse the tab-separated percentile dump down to
# a single value (matches the corrected README stats).
$t.Cell(44,1).Range.Text = "99.41"
$t.Cell(45,1).Range.Text = "1.3"
$t.Cell(46,1).Range.Text = "220"
